$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 100,2
$data[0,0] = 498
$data[0,1] = 289
$data[1,0] = 508
$data[1,1] = 280
$data[2,0] = 472
$data[2,1] = 210
$data[3,0] = 420
$data[3,1] = 290
$data[4,0] = 518
$data[4,1] = 230
$data[5,0] = 507
$data[5,1] = 279
$data[6,0] = 419
$data[6,1] = 289
$data[7,0] = 466
$data[7,1] = 279
$data[8,0] = 407
$data[8,1] = 190
$data[9,0] = 407
$data[9,1] = 190
$data[10,0] = 528
$data[10,1] = 220
$data[11,0] = 466
$data[11,1] = 279
$data[12,0] = 507
$data[12,1] = 279
$data[13,0] = 407
$data[13,1] = 190
$data[14,0] = 407
$data[14,1] = 190
$data[15,0] = 407
$data[15,1] = 190
$data[16,0] = 466
$data[16,1] = 279
$data[17,0] = 407
$data[17,1] = 190
$data[18,0] = 419
$data[18,1] = 289
$data[19,0] = 407
$data[19,1] = 190
$data[20,0] = 407
$data[20,1] = 190
$data[21,0] = 407
$data[21,1] = 190
$data[22,0] = 442
$data[22,1] = 279
$data[23,0] = 407
$data[23,1] = 190
$data[24,0] = 457
$data[24,1] = 289
$data[25,0] = 407
$data[25,1] = 190
$data[26,0] = 407
$data[26,1] = 190
$data[27,0] = 518
$data[27,1] = 230
$data[28,0] = 407
$data[28,1] = 190
$data[29,0] = 517
$data[29,1] = 229
$data[30,0] = 431
$data[30,1] = 290
$data[31,0] = 435
$data[31,1] = 279
$data[32,0] = 507
$data[32,1] = 279
$data[33,0] = 458
$data[33,1] = 290
$data[34,0] = 407
$data[34,1] = 190
$data[35,0] = 431
$data[35,1] = 290
$data[36,0] = 436
$data[36,1] = 280
$data[37,0] = 498
$data[37,1] = 289
$data[38,0] = 407
$data[38,1] = 190
$data[39,0] = 407
$data[39,1] = 190
$data[40,0] = 407
$data[40,1] = 190
$data[41,0] = 402
$data[41,1] = 190
$data[42,0] = 523
$data[42,1] = 220
$data[43,0] = 402
$data[43,1] = 190
$data[44,0] = 402
$data[44,1] = 190
$data[45,0] = 402
$data[45,1] = 190
$data[46,0] = 493
$data[46,1] = 289
$data[47,0] = 503
$data[47,1] = 280
$data[48,0] = 435
$data[48,1] = 289
$data[49,0] = 402
$data[49,1] = 190
$data[50,0] = 494
$data[50,1] = 290
$data[51,0] = 412
$data[51,1] = 289
$data[52,0] = 448
$data[52,1] = 290
$data[53,0] = 431
$data[53,1] = 280
$data[54,0] = 402
$data[54,1] = 190
$data[55,0] = 402
$data[55,1] = 190
$data[56,0] = 401
$data[56,1] = 189
$data[57,0] = 402
$data[57,1] = 190
$data[58,0] = 502
$data[58,1] = 279
$data[59,0] = 447
$data[59,1] = 289
$data[60,0] = 402
$data[60,1] = 190
$data[61,0] = 401
$data[61,1] = 189
$data[62,0] = 402
$data[62,1] = 190
$data[63,0] = 495
$data[63,1] = 209
$data[64,0] = 503
$data[64,1] = 280
$data[65,0] = 402
$data[65,1] = 190
$data[66,0] = 503
$data[66,1] = 280
$data[67,0] = 431
$data[67,1] = 280
$data[68,0] = 402
$data[68,1] = 190
$data[69,0] = 402
$data[69,1] = 190
$data[70,0] = 447
$data[70,1] = 289
$data[71,0] = 402
$data[71,1] = 190
$data[72,0] = 523
$data[72,1] = 220
$data[73,0] = 447
$data[73,1] = 289
$data[74,0] = 451
$data[74,1] = 250
$data[75,0] = 402
$data[75,1] = 190
$data[76,0] = 402
$data[76,1] = 190
$data[77,0] = 523
$data[77,1] = 220
$data[78,0] = 402
$data[78,1] = 190
$data[79,0] = 402
$data[79,1] = 190
$data[80,0] = 402
$data[80,1] = 190
$data[81,0] = 451
$data[81,1] = 250
$data[82,0] = 402
$data[82,1] = 190
$data[83,0] = 432
$data[83,1] = 289
$data[84,0] = 402
$data[84,1] = 190
$data[85,0] = 412
$data[85,1] = 289
$data[86,0] = 405
$data[86,1] = 289
$data[87,0] = 451
$data[87,1] = 250
$data[88,0] = 401
$data[88,1] = 189
$data[89,0] = 402
$data[89,1] = 190
$data[90,0] = 402
$data[90,1] = 190
$data[91,0] = 461
$data[91,1] = 279
$data[92,0] = 402
$data[92,1] = 190
$data[93,0] = 496
$data[93,1] = 210
$data[94,0] = 402
$data[94,1] = 190
$data[95,0] = 467
$data[95,1] = 210
$data[96,0] = 402
$data[96,1] = 190
$data[97,0] = 402
$data[97,1] = 190
$data[98,0] = 503
$data[98,1] = 280
$data[99,0] = 402
$data[99,1] = 190

$range = $ws.Range("B2:C101")
$range.Value = $data
